# Feb 12 update 3
# Update Sheet1: revise the budget-category answer text in B7 (second
# stimulus "spend on" answers), align the used range to top-vertical,
# fix up row 9's autofit height, and leave the cursor on B10 (matches
# the author's last edit position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) B7's answer text was retyped - same bullet list, but the
#    "Mortgage ... Utilities" line lost two stray spaces, so Excel treats
#    it as new, distinct text (new shared string) rather than reusing the
#    untouched string that B4/B10 still point to.
$ws.Range("B7").Value = "•Food`n•Clothing`n•Household supplies and personal care                                 •Household items (TV, electronics, furniture, appliances)  `n•Recreational goods (sports and fitness equipment, bicycles, toys, games)                             `n•Rent`n•Mortgage                                                                             •Utilities and telecommunications`n•Vehicle payments                                                              •Paying down credit card, student loans, or other debts  `n•Charitable donations or giving to family members                             `n•Savings or other investments`n•Other, please describe:"

# 2) Select-all + Align Top across every populated cell.
$topAlign = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop
$ws.Range("A1:E10").VerticalAlignment = $topAlign
$ws.Range("A11:A16").VerticalAlignment = $topAlign

# 3) Row 9 re-wrapped very slightly shorter once the alignment/content
#    settled (238 -> unaffected, but row 9 specifically: 106 -> 105).
$ws.Rows.Item(9).RowHeight = 105

# 4) Last selection before save was B10 (not A10).
$ws.Range("B10").Select()

Write-Host "edit applied"
